# This script applies a "re-sort" of rows 59-77 on the active worksheet.
# The row-level metadata (location, county, dates, reporter, etc.) stays the
# same per physical row, but the observation-specific fields are permuted
# among the 19 rows according to the mapping derived from the source diff.
# Additionally the Ost/Nord (Q/R) coordinate columns, which previously held
# long-decimal values, are rounded to whole numbers in the process.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose VALUES move between rows (everything else on the row is
# left untouched).
$cols = @("A","B","D","E","F","G","H","I","M","Q","R","S","Z")

$firstRow = 59
$lastRow = 77

# Snapshot the current contents of the moving columns for every row in the
# block, keyed by row number.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# target row -> source row (which row's data ends up at this row)
$mapping = @{
    59 = 61
    60 = 73
    61 = 75
    62 = 63
    63 = 64
    64 = 62
    65 = 67
    66 = 74
    67 = 72
    68 = 65
    69 = 68
    70 = 59
    71 = 77
    72 = 76
    73 = 70
    74 = 71
    75 = 66
    76 = 60
    77 = 69
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $src = $snapshot[$sourceRow]

    foreach ($c in $cols) {
        $val = $src[$c]

        if ($c -eq "Q" -or $c -eq "R") {
            $val = [math]::Floor([double]$val + 0.5)
        }

        $ws.Range("$c$targetRow").Value = $val
    }

    # Slutttid (AB) always mirrors Starttid (Z)
    $ws.Range("AB$targetRow").Value = $src["Z"]
}
